# Applies the row-content swap described in the commit diff.
# Pairs of rows (10/12, 13/15, 20/21, 23/26, 24/25, 27/28) exchange their
# observation-identity fields (species, counts, coordinates, external id,
# times, observer) while the shared locality/date columns stay put.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("A10").Value = 131106319
$ws.Range("B10").Value = 92107
$ws.Range("D10").Value = 'NT'
$ws.Range("E10").Value = 658
$ws.Range("F10").Value = 'Rosenticka'
$ws.Range("G10").Value = 'Fomitopsis rosea'
$ws.Range("H10").Value = '(Alb. & Schwein.:Fr.) P.Karst.'
$ws.Range("I10").Value = '2'
$ws.Range("Q10").Value = 601569
$ws.Range("R10").Value = 6992657
$ws.Range("X10").Value = '2025_0864'
$ws.Range("Z10").Value = '13:14'
$ws.Range("AB10").Value = '13:14'

# Row 12
$ws.Range("A12").Value = 131106321
$ws.Range("B12").Value = 92022
$ws.Range("D12").Value = 'VU'
$ws.Range("E12").Value = 48
$ws.Range("F12").Value = 'Lappticka'
$ws.Range("G12").Value = 'Amylocystis lapponica'
$ws.Range("H12").Value = '(Romell) Bondartsev & Singer'
$ws.Range("I12").Value = ""
$ws.Range("Q12").Value = 601579
$ws.Range("R12").Value = 6992698
$ws.Range("X12").Value = '2025_0862'
$ws.Range("Z12").Value = '12:56'
$ws.Range("AB12").Value = '12:56'

# Row 13
$ws.Range("A13").Value = 131106314
$ws.Range("I13").Value = '1'
$ws.Range("Q13").Value = 601556
$ws.Range("R13").Value = 6992605
$ws.Range("X13").Value = '2025_0870'
$ws.Range("Z13").Value = '13:21'
$ws.Range("AB13").Value = '13:21'
$ws.Range("AX13").Value = 'David Isaksson'

# Row 15
$ws.Range("A15").Value = 131106325
$ws.Range("I15").Value = ""
$ws.Range("Q15").Value = 601615
$ws.Range("R15").Value = 6992785
$ws.Range("X15").Value = '2025_0858'
$ws.Range("Z15").Value = '12:21'
$ws.Range("AB15").Value = '12:21'
$ws.Range("AX15").Value = 'Alexander Hoffmann'

# Row 20
$ws.Range("A20").Value = 131106310
$ws.Range("I20").Value = ""
$ws.Range("J20").Value = ""
$ws.Range("Q20").Value = 601470
$ws.Range("R20").Value = 6992568
$ws.Range("X20").Value = '2025_0874'
$ws.Range("Z20").Value = '13:35'
$ws.Range("AB20").Value = '13:35'
$ws.Range("AX20").Value = 'Alexander Hoffmann'

# Row 21
$ws.Range("A21").Value = 131106313
$ws.Range("I21").Value = '1'
$ws.Range("J21").Value = 'mycel'
$ws.Range("Q21").Value = 601530
$ws.Range("R21").Value = 6992589
$ws.Range("X21").Value = '2025_0871'
$ws.Range("Z21").Value = '13:29'
$ws.Range("AB21").Value = '13:29'
$ws.Range("AX21").Value = 'David Isaksson'

# Row 23
$ws.Range("A23").Value = 131106311
$ws.Range("B23").Value = 91809
$ws.Range("E23").Value = 1202
$ws.Range("F23").Value = 'Ullticka'
$ws.Range("G23").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H23").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("I23").Value = ""
$ws.Range("J23").Value = ""
$ws.Range("Q23").Value = 601498
$ws.Range("R23").Value = 6992583
$ws.Range("X23").Value = '2025_0873'
$ws.Range("Z23").Value = '13:32'
$ws.Range("AB23").Value = '13:32'
$ws.Range("AC23").Value = ""
$ws.Range("AX23").Value = 'Alexander Hoffmann'

# Row 24
$ws.Range("A24").Value = 131106327
$ws.Range("B24").Value = 91809
$ws.Range("E24").Value = 1202
$ws.Range("F24").Value = 'Ullticka'
$ws.Range("G24").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H24").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("J24").Value = 'mycel'
$ws.Range("R24").Value = 6992789
$ws.Range("X24").Value = '2025_0856'
$ws.Range("Z24").Value = '12:10'
$ws.Range("AB24").Value = '12:10'
$ws.Range("AX24").Value = 'David Isaksson, Alexander Hoffmann'

# Row 25
$ws.Range("A25").Value = 131106323
$ws.Range("B25").Value = 92107
$ws.Range("E25").Value = 658
$ws.Range("F25").Value = 'Rosenticka'
$ws.Range("G25").Value = 'Fomitopsis rosea'
$ws.Range("H25").Value = '(Alb. & Schwein.:Fr.) P.Karst.'
$ws.Range("J25").Value = ""
$ws.Range("R25").Value = 6992738
$ws.Range("X25").Value = '2025_0860'
$ws.Range("Z25").Value = '12:35'
$ws.Range("AB25").Value = '12:35'
$ws.Range("AX25").Value = 'David Isaksson'

# Row 26
$ws.Range("A26").Value = 131108353
$ws.Range("B26").Value = 10966
$ws.Range("E26").Value = 101449
$ws.Range("F26").Value = ""
$ws.Range("G26").Value = 'Olisthaerus substriatus'
$ws.Range("H26").Value = '(Paykull, 1790)'
$ws.Range("I26").Value = '1'
$ws.Range("J26").Value = 'ex.'
$ws.Range("Q26").Value = 601612
$ws.Range("R26").Value = 6992796
$ws.Range("X26").Value = '2025_0855'
$ws.Range("Z26").Value = '12:09'
$ws.Range("AB26").Value = '12:09'
$ws.Range("AC26").Value = 'Granlåga med både rynkskinn, ull-och violticka'
$ws.Range("AX26").Value = 'David Isaksson'

# Row 27
$ws.Range("A27").Value = 131106330
$ws.Range("B27").Value = 92107
$ws.Range("D27").Value = 'NT'
$ws.Range("E27").Value = 658
$ws.Range("F27").Value = 'Rosenticka'
$ws.Range("G27").Value = 'Fomitopsis rosea'
$ws.Range("H27").Value = '(Alb. & Schwein.:Fr.) P.Karst.'
$ws.Range("J27").Value = ""
$ws.Range("Q27").Value = 601607
$ws.Range("R27").Value = 6992782
$ws.Range("X27").Value = '2025_0853'
$ws.Range("Z27").Value = '12:06'
$ws.Range("AB27").Value = '12:06'

# Row 28
$ws.Range("A28").Value = 131106329
$ws.Range("B28").Value = 92268
$ws.Range("D28").Value = 'VU'
$ws.Range("E28").Value = 1209
$ws.Range("F28").Value = 'Rynkskinn'
$ws.Range("G28").Value = 'Hermanssonia centrifuga'
$ws.Range("H28").Value = '(P. Karst.) Zmitr.'
$ws.Range("J28").Value = 'mycel'
$ws.Range("Q28").Value = 601609
$ws.Range("R28").Value = 6992789
$ws.Range("X28").Value = '2025_0854'
$ws.Range("Z28").Value = '12:09'
$ws.Range("AB28").Value = '12:09'
